# Replace the "Zippers" bullet text with its reworded / reordered version.
#
# Before: "Aggregation / Addressing: Locations / Contexts, Streams, Reactive.
#          Monads, Augmentation (navigation / transforms). Zippers."
# After:  "Zippers: Aggregation / Addressing: Locations / Contexts. Parsing.
#          Monads. Augmentation (navigation / transforms) Reactive Streams."

$d = $word.ActiveDocument

$old = "Aggregation / Addressing: Locations / Contexts, Streams, Reactive. Monads, Augmentation (navigation / transforms). Zippers."
$new = "Zippers: Aggregation / Addressing: Locations / Contexts. Parsing. Monads. Augmentation (navigation / transforms) Reactive Streams."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)
